# Update gh-pages to output generated at 456a3b4
# Apply refreshed "想去人数" (F) / "最低票价" (G) values across sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = 70
$ws1.Range("F5").Value = 241
$ws1.Range("F7").Value = 181
$ws1.Range("F8").Value = 3
$ws1.Range("F9").Value = 6380
$ws1.Range("F13").Value = 5564
$ws1.Range("F16").Value = 1216
$ws1.Range("F18").Value = 70
$ws1.Range("F20").Value = 78
$ws1.Range("F25").Value = 3991
$ws1.Range("F26").Value = 13
$ws1.Range("F27").Value = 177

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 103

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 103
$ws4.Range("G4").Value = 70
$ws4.Range("F6").Value = 241
$ws4.Range("F8").Value = 181
$ws4.Range("F9").Value = 3
$ws4.Range("F10").Value = 6380
$ws4.Range("F14").Value = 5564
$ws4.Range("F17").Value = 1216
$ws4.Range("F19").Value = 70
$ws4.Range("F21").Value = 78
$ws4.Range("F26").Value = 3991
$ws4.Range("F28").Value = 13
$ws4.Range("F29").Value = 177
